# Weekly update: insert a new data row at row 8 (pushing existing rows 8-46
# down to 9-47) and populate it with this week's "Haba" price record for the
# Terminal La Palmera de La Serena market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8:46 down to 9:47, creating a fresh blank row 8.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the latest weekly record.
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 45050
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100112026
$ws.Cells.Item(8, 7).Value = "Haba"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 360
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14500
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 580
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
